$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = '43.003.99'
$ws.Range("E2").Value2 = '  +0.84%  '
$ws.Range("D3").Value2 = '2.543.05'
$ws.Range("E3").Value2 = '  +0.61%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = '0.998'
$ws.Range("E4").Value2 = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '316.95'
$ws.Range("E5").Value2 = '  +0.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '97.63'
$ws.Range("E6").Value2 = '  +2.96%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.575'
$ws.Range("E7").Value2 = '  -0.63%  '
$ws.Range("E8").Value2 = '  -0.11%  '
$ws.Range("E9").Value2 = '  -0.45%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '36.33'
$ws.Range("E10").Value2 = '  +0.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '0.0821'
$ws.Range("E11").Value2 = '  +1.21%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '7.65'
$ws.Range("E12").Value2 = '  -0.55%  '
$ws.Range("E13").Value2 = '  -2.96%  '
$ws.Range("D14").Value2 = '2.929.90'
$ws.Range("E14").Value2 = '  +0.59%  '
$ws.Range("D15").Value2 = '2.546.17'
$ws.Range("E15").Value2 = '  +1.66%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '15.21'
$ws.Range("E16").Value2 = '  -2.06%  '
$ws.Range("E17").Value2 = '  -1.25%  '
$ws.Range("D18").Value2 = '42.993.00'
$ws.Range("E18").Value2 = '  +0.61%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '6.87'
$ws.Range("E19").Value2 = '  +5.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '12.83'
$ws.Range("E20").Value2 = '  -1.80%  '
$ws.Range("E21").Value2 = '  +0.25%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '69.99'
$ws.Range("E22").Value2 = '  -1.77%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '254.69'
$ws.Range("E23").Value2 = '  +0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '2.95'
$ws.Range("E24").Value2 = '  -1.23%  '
$ws.Range("E25").Value2 = '  +0.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '26.63'
$ws.Range("E26").Value2 = '  -3.75%  '
$ws.Range("E27").Value2 = '  +0.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '2.42'
$ws.Range("E28").Value2 = '  +5.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '41.07'
$ws.Range("E29").Value2 = '  +4.44%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '10.47'
$ws.Range("E30").Value2 = '  +4.14%  '
$ws.Range("E31").Value2 = '  +0.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '158.38'
$ws.Range("E32").Value2 = '  +1.35%  '
$ws.Range("E33").Value2 = '  +3.67%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '3.37'
$ws.Range("E34").Value2 = '  +0.93%  '
$ws.Range("E35").Value2 = '  +4.62%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '19.04'
$ws.Range("E36").Value2 = '  -4.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '0.0795'
$ws.Range("E37").Value2 = '  +1.23%  '
$ws.Range("E38").Value2 = '  +0.34%  '
$ws.Range("E39").Value2 = '  +14.62%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '0.119'
$ws.Range("E40").Value2 = '  -0.60%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '21.92'
$ws.Range("E41").Value2 = '  -11.68%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '3.85'
$ws.Range("E42").Value2 = '  -0.11%  '
$ws.Range("E43").Value2 = '  +0.28%  '
$ws.Range("E44").Value2 = '  +0.15%  '
$ws.Range("E45").Value2 = '  -2.17%  '
$ws.Range("D46").Value2 = '2.025.89'
$ws.Range("E46").Value2 = '  -1.31%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '9.12'
$ws.Range("E47").Value2 = '  +3.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '84.64'
$ws.Range("E48").Value2 = '  -2.13%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '77.12'
$ws.Range("E49").Value2 = '  +3.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '106.75'
$ws.Range("E50").Value2 = '  +4.86%  '
$ws.Range("D51").Value2 = '2.785.05'
$ws.Range("E51").Value2 = '  +0.88%  '
